$d = $word.ActiveDocument

$replacements = @(
    @{old = "41×60=2460"; new = "41×62=2542"},
    @{old = "41×89=3649"; new = "56×16=896"},
    @{old = "77×11=847"; new = "51×19=969"},
    @{old = "48×81=3888"; new = "19×28=532"},
    @{old = "21×81=1701"; new = "93×63=5859"},
    @{old = "68×57=3876"; new = "81×93=7533"},
    @{old = "94×32=3008"; new = "95×40=3800"},
    @{old = "77×69=5313"; new = "62×93=5766"},
    @{old = "21×58=1218"; new = "49×13=637"},
    @{old = "74×22=1628"; new = "77×43=3311"},
    @{old = "45×95=4275"; new = "49×79=3871"},
    @{old = "47×78=3666"; new = "15×94=1410"},
    @{old = "11×74=814"; new = "85×15=1275"},
    @{old = "89×18=1602"; new = "43×83=3569"},
    @{old = "61×67=4087"; new = "94×20=1880"},
    @{old = "73×19=1387"; new = "13×82=1066"},
    @{old = "67×88=5896"; new = "81×44=3564"},
    @{old = "17×71=1207"; new = "27×38=1026"},
    @{old = "36×96=3456"; new = "90×96=8640"},
    @{old = "28×39=1092"; new = "40×99=3960"},
    @{old = "78×20=1560"; new = "78×52=4056"},
    @{old = "53×66=3498"; new = "95×81=7695"},
    @{old = "33×95=3135"; new = "39×63=2457"},
    @{old = "41×77=3157"; new = "93×48=4464"},
    @{old = "62×58=3596"; new = "54×99=5346"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $r.new, 2)
}
